$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 5
$ws.Range("B2").Value = 6
$ws.Range("B3").Value = 7
$ws.Range("B4").Value = 8

$ws.Range("H11").Select()
